$d = $word.ActiveDocument
$newText = "V roku Gemini: 14.-23. februára, 14.-24. marca"

# Walk every paragraph looking for the old "Perseus" date blurb (it may be
# split across several runs, e.g. "...súhvezdie " + "Perseus" + ": " + "30. ...").
# Replace the whole paragraph's text with a single, plain run (no rPr) that
# holds the translated Gemini dates, matching what Word does when a user
# selects the run(s), deletes them, and types fresh text into the gap.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -like "*V roku 2018*") {
        $r = $p.Range
        # Exclude the trailing paragraph mark from the range to delete.
        $r.MoveEnd(1, -1) | Out-Null
        $r.Delete() | Out-Null
        $r.InsertAfter($newText) | Out-Null
    }
}
